$wb = $excel.ActiveWorkbook

$partidos = $wb.Worksheets.Item(1)
$atributos = $wb.Worksheets.Item(2)

# --- Update selection on "partidos" (sheet1) ---
$partidos.Range("H40").Select()

# --- Update selection on "atributos" (sheet2) ---
$atributos.Range("A1:A24").Select()

# --- Add the new "analisis" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$analisis = $wb.Worksheets.Add($null, $lastSheet)
$analisis.Name = "analisis"

# --- Header row ---
$analisis.Range("A1").Value = "Jugador"
$analisis.Range("B1").Value = "Numero de partidos"
$analisis.Range("C1").Value = "Ratio de victorias"

$atributos.Range("A1").Copy()
$analisis.Range("A1:C1").PasteSpecial(-4122)

# --- Player rows ---
$players = @("CMC","Moi","Lolo","Celia","Ana Lucia","Paco Ch","Migue","Pabliyo","Pau","Héctor","Javi","Sonia","Coca","Carlos López","Diego","Jesus","Jesus María","Paco Cádiz","Richard","Juanma","Dani","Martin","El Largo")

for ($i = 0; $i -lt $players.Count; $i++) {
    $row = $i + 2
    $name = $players[$i]

    $aCell = $analisis.Cells.Item($row, 1)
    $bCell = $analisis.Cells.Item($row, 2)
    $cCell = $analisis.Cells.Item($row, 3)

    $aCell.Value = $name

    $atributos.Range("A2").Copy()
    $aCell.PasteSpecial(-4122)

    $bCell.Formula = '=COUNTIF(partidos!$C$2:$C$363,A' + $row + ')'
    $cCell.Formula = '=COUNTIFS(partidos!C:C, A' + $row + ', partidos!E:E, "Gana")/B' + $row
}

# --- Column widths for the new sheet ---
$analisis.Columns.Item(1).ColumnWidth = 17.166666666666668
$analisis.Columns.Item(2).ColumnWidth = 23.166666666666668
$analisis.Columns.Item(3).ColumnWidth = 17.451822916666668

# --- Row height for the "Carlos López" row (row 15) ---
$analisis.Rows.Item(15).RowHeight = 16.5

# --- Final selection + make this the active sheet/tab ---
$analisis.Range("D3").Select()
